$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows 27-29 (columns A, B, E, F, G, H) are cyclically shifted:
# old row27 values -> row29, old row28 values -> row27, old row29 values -> row28

# Capture current (before) values for the affected columns
$A27 = $ws.Range("A27").Value()
$B27 = $ws.Range("B27").Value()
$E27 = $ws.Range("E27").Value()
$F27 = $ws.Range("F27").Value()
$G27 = $ws.Range("G27").Value()
$H27 = $ws.Range("H27").Value()

$A28 = $ws.Range("A28").Value()
$B28 = $ws.Range("B28").Value()
$E28 = $ws.Range("E28").Value()
$F28 = $ws.Range("F28").Value()
$G28 = $ws.Range("G28").Value()
$H28 = $ws.Range("H28").Value()

$A29 = $ws.Range("A29").Value()
$B29 = $ws.Range("B29").Value()
$E29 = $ws.Range("E29").Value()
$F29 = $ws.Range("F29").Value()
$G29 = $ws.Range("G29").Value()
$H29 = $ws.Range("H29").Value()

# Row 27 gets old row 28 values
$ws.Range("A27").Value = $A28
$ws.Range("B27").Value = $B28
$ws.Range("E27").Value = $E28
$ws.Range("F27").Value = $F28
$ws.Range("G27").Value = $G28
$ws.Range("H27").Value = $H28

# Row 28 gets old row 29 values
$ws.Range("A28").Value = $A29
$ws.Range("B28").Value = $B29
$ws.Range("E28").Value = $E29
$ws.Range("F28").Value = $F29
$ws.Range("G28").Value = $G29
$ws.Range("H28").Value = $H29

# Row 29 gets old row 27 values
$ws.Range("A29").Value = $A27
$ws.Range("B29").Value = $B27
$ws.Range("E29").Value = $E27
$ws.Range("F29").Value = $F27
$ws.Range("G29").Value = $G27
$ws.Range("H29").Value = $H27
